$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'305.31"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'0.08%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'35.75"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'-0.44%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'5.041"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'-1.02%"
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'0.08031"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'-0.39%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'1.873"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'-2.68%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = "'4.141"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'-0.80%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').Value = "'7.794"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'0.60%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = "'0.9198"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'-0.90%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = "'0.1273"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'-6.38%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = "'0.1913"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'0.57%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = "'0.09086"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'-0.77%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.03472"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'-4.65%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09863"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'0.46%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001408"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'-0.57%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = "'0.006213"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'5.23%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = "'3.821"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'7.32%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('D19').Value = "'0.3419"
$ws.Range('D19').Style = 'Normal'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').Value = "'0.1320"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'-0.90%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('B21').Value = 'MCDex'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D21').Value = "'5.217"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'6.57%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.2305"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'-11.41%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.04428"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'-0.46%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'0.001234"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'0.99%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.004613"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'-3.86%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('E27').Value = "'-3.82%"
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'42.05%"
$ws.Range('E28').Style = 'Normal'
$ws.Range('D39').Value = "'0.01946"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'-1.24%"
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.05290"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'7.51%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.007609"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'-0.42%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'10.64%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'-1.61%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.002162"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'2.90%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'-15.21%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.00006137"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'-3.65%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'0.03%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'63.63"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'0.10%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'39.37%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'0.03%"
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'0.03%"
$ws.Range('E51').Style = 'Normal'
